# Weekly update for "Hortaliza, Comercializadora del Agro de Limarí - Poroto granado":
# a new observation is inserted as row 51 (pushing the existing rows 51-73 down to
# 52-74, i.e. one new week of data is prepended to the series) and the sheet
# dimension grows from A1:R73 to A1:R74.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 51; Excel shifts rows 51:73 down to 52:74 and extends the
# used range / dimension automatically.
$ws.Rows("51:51").Insert()

# The columns that are constant for every record in this sheet (market id,
# market name, region, codreg, category id/name, variety, quality, origin,
# classification) simply need to be copied down from the row directly below
# (which used to be row 51 before the insert pushed it to row 52).
$ws.Range("A51").Value = $ws.Range("A52").Value2
$ws.Range("B51").Value = $ws.Range("B52").Value2
$ws.Range("C51").Value = $ws.Range("C52").Value2
$ws.Range("E51").Value = $ws.Range("E52").Value2
$ws.Range("F51").Value = $ws.Range("F52").Value2
$ws.Range("G51").Value = $ws.Range("G52").Value2
$ws.Range("H51").Value = $ws.Range("H52").Value2
$ws.Range("I51").Value = $ws.Range("I52").Value2
$ws.Range("O51").Value = $ws.Range("O52").Value2
$ws.Range("R51").Value = $ws.Range("R52").Value2

# New weekly record's own values (date, volume, min/max/weighted-avg price,
# unit of sale and its price-per-kg / kg-equivalent breakdown).
$ws.Range("D51").Value = 44609
$ws.Range("J51").Value = 600
$ws.Range("K51").Value = 22000
$ws.Range("L51").Value = 24000
$ws.Range("M51").Value = 23000
$ws.Range("N51").Value = "$/malla 25 kilos"
$ws.Range("P51").Value = 920
$ws.Range("Q51").Value = 25
